$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.440918333333333
$ws.Range("H2").Value = 7.322755
$ws.Range("I2").Value = 0.5182826554654038
$ws.Range("J2").Value = 0.5182826554654038
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 163.7119853333333
$ws.Range("N2").Value = 491.135956
$ws.Range("O2").Value = 0.2754003062401033
$ws.Range("P2").Value = 0.2754003062401033
$ws.Range("Q2").Value = 399.6075863865311
$ws.Range("R2").Value = 3596.46827747878
$ws.Range("S2").Value = 0.1427352020341061
$ws.Range("T2").Value = 0.1427352020341061

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.440918333333333
$ws.Range("H3").Value = 7.322755
$ws.Range("I3").Value = 0.5182826554654038
$ws.Range("J3").Value = 0.5182826554654038
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.290282124557779
$ws.Range("P3").Value = 0.290282124557779
$ws.Range("Q3").Value = 421.201198899745
$ws.Range("R3").Value = 3790.810790097705
$ws.Range("S3").Value = 0.1504481903499448
$ws.Range("T3").Value = 0.1504481903499448

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.440918333333333
$ws.Range("H4").Value = 7.322755
$ws.Range("I4").Value = 0.5182826554654038
$ws.Range("J4").Value = 0.5182826554654038
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 111.4881643333333
$ws.Range("N4").Value = 334.464493
$ws.Range("O4").Value = 0.1875481171218523
$ws.Range("P4").Value = 0.1875481171218523
$ws.Range("Q4").Value = 272.1335042709127
$ws.Range("R4").Value = 2449.201538438215
$ws.Range("S4").Value = 0.09720293616945017
$ws.Range("T4").Value = 0.09720293616945018

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.440918333333333
$ws.Range("H5").Value = 7.322755
$ws.Range("I5").Value = 0.5182826554654038
$ws.Range("J5").Value = 0.5182826554654038
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 86.95798233333333
$ws.Range("N5").Value = 260.873947
$ws.Range("O5").Value = 0.1462828449356383
$ws.Range("P5").Value = 0.1462828449356383
$ws.Range("Q5").Value = 212.2573333071094
$ws.Range("R5").Value = 1910.315999763985
$ws.Range("S5").Value = 0.07581586132227652
$ws.Range("T5").Value = 0.07581586132227652

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.440918333333333
$ws.Range("H6").Value = 7.322755
$ws.Range("I6").Value = 0.5182826554654038
$ws.Range("J6").Value = 0.5182826554654038
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 59.73436333333333
$ws.Range("N6").Value = 179.20309
$ws.Range("O6").Value = 0.100486607144627
$ws.Range("P6").Value = 0.100486607144627
$ws.Range("Q6").Value = 145.8067025903278
$ws.Range("R6").Value = 1312.26032331295
$ws.Range("S6").Value = 0.05208046558962608
$ws.Range("T6").Value = 0.05208046558962608

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6061233333333333
$ws.Range("H7").Value = 1.81837
$ws.Range("I7").Value = 0.1286987796558298
$ws.Range("J7").Value = 0.1286987796558298
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 163.7119853333333
$ws.Range("N7").Value = 491.135956
$ws.Range("O7").Value = 0.2754003062401033
$ws.Range("P7").Value = 0.2754003062401033
$ws.Range("Q7").Value = 99.22965425685778
$ws.Range("R7").Value = 893.0668883117199
$ws.Range("S7").Value = 0.03544368332994311
$ws.Range("T7").Value = 0.03544368332994312

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6061233333333333
$ws.Range("H8").Value = 1.81837
$ws.Range("I8").Value = 0.1286987796558298
$ws.Range("J8").Value = 0.1286987796558298
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.290282124557779
$ws.Range("P8").Value = 0.290282124557779
$ws.Range("Q8").Value = 104.59173139663
$ws.Range("R8").Value = 941.3255825696699
$ws.Range("S8").Value = 0.03735895518648774
$ws.Range("T8").Value = 0.03735895518648775

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6061233333333333
$ws.Range("H9").Value = 1.81837
$ws.Range("I9").Value = 0.1286987796558298
$ws.Range("J9").Value = 0.1286987796558298
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 111.4881643333333
$ws.Range("N9").Value = 334.464493
$ws.Range("O9").Value = 0.1875481171218523
$ws.Range("P9").Value = 0.1875481171218523
$ws.Range("Q9").Value = 67.57557779293444
$ws.Range("R9").Value = 608.18020013641
$ws.Range("S9").Value = 0.02413721380033103
$ws.Range("T9").Value = 0.02413721380033104

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6061233333333333
$ws.Range("H10").Value = 1.81837
$ws.Range("I10").Value = 0.1286987796558298
$ws.Range("J10").Value = 0.1286987796558298
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 86.95798233333333
$ws.Range("N10").Value = 260.873947
$ws.Range("O10").Value = 0.1462828449356383
$ws.Range("P10").Value = 0.1462828449356383
$ws.Range("Q10").Value = 52.70726211182112
$ws.Range("R10").Value = 474.36535900639
$ws.Range("S10").Value = 0.01882642362779964
$ws.Range("T10").Value = 0.01882642362779964

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.6061233333333333
$ws.Range("H11").Value = 1.81837
$ws.Range("I11").Value = 0.1286987796558298
$ws.Range("J11").Value = 0.1286987796558298
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 59.73436333333333
$ws.Range("N11").Value = 179.20309
$ws.Range("O11").Value = 0.100486607144627
$ws.Range("P11").Value = 0.100486607144627
$ws.Range("Q11").Value = 36.20639141814444
$ws.Range("R11").Value = 325.8575227633
$ws.Range("S11").Value = 0.01293250371126828
$ws.Range("T11").Value = 0.01293250371126828

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.662586
$ws.Range("H12").Value = 4.987757999999999
$ws.Range("I12").Value = 0.3530185648787664
$ws.Range("J12").Value = 0.3530185648787664
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 163.7119853333333
$ws.Range("N12").Value = 491.135956
$ws.Range("O12").Value = 0.2754003062401033
$ws.Range("P12").Value = 0.2754003062401033
$ws.Range("Q12").Value = 272.1852548474053
$ws.Range("R12").Value = 2449.667293626648
$ws.Range("S12").Value = 0.09722142087605404
$ws.Range("T12").Value = 0.09722142087605404

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.662586
$ws.Range("H13").Value = 4.987757999999999
$ws.Range("I13").Value = 0.3530185648787664
$ws.Range("J13").Value = 0.3530185648787664
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.290282124557779
$ws.Range("P13").Value = 0.290282124557779
$ws.Range("Q13").Value = 286.893341293242
$ws.Range("R13").Value = 2582.040071639178
$ws.Range("S13").Value = 0.1024749790213464
$ws.Range("T13").Value = 0.1024749790213464

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.662586
$ws.Range("H14").Value = 4.987757999999999
$ws.Range("I14").Value = 0.3530185648787664
$ws.Range("J14").Value = 0.3530185648787664
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 111.4881643333333
$ws.Range("N14").Value = 334.464493
$ws.Range("O14").Value = 0.1875481171218523
$ws.Range("P14").Value = 0.1875481171218523
$ws.Range("Q14").Value = 185.3586611862993
$ws.Range("R14").Value = 1668.227950676694
$ws.Range("S14").Value = 0.06620796715207108
$ws.Range("T14").Value = 0.0662079671520711

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.662586
$ws.Range("H15").Value = 4.987757999999999
$ws.Range("I15").Value = 0.3530185648787664
$ws.Range("J15").Value = 0.3530185648787664
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 86.95798233333333
$ws.Range("N15").Value = 260.873947
$ws.Range("O15").Value = 0.1462828449356383
$ws.Range("P15").Value = 0.1462828449356383
$ws.Range("Q15").Value = 144.5751240156473
$ws.Range("R15").Value = 1301.176116140826
$ws.Range("S15").Value = 0.05164055998556216
$ws.Range("T15").Value = 0.05164055998556216

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.662586
$ws.Range("H16").Value = 4.987757999999999
$ws.Range("I16").Value = 0.3530185648787664
$ws.Range("J16").Value = 0.3530185648787664
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 59.73436333333333
$ws.Range("N16").Value = 179.20309
$ws.Range("O16").Value = 0.100486607144627
$ws.Range("P16").Value = 0.100486607144627
$ws.Range("Q16").Value = 99.31351619691331
$ws.Range("R16").Value = 893.8216457722198
$ws.Range("S16").Value = 0.03547363784373261
$ws.Range("T16").Value = 0.03547363784373261
